$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - LinearRegression (values updated slightly)
$ws.Range("B2").Value = -69.02463949302539
$ws.Range("C2").Value = -69.02463949302539
$ws.Range("D2").Value = -69.02463949302539

# Row 3 - RandomForestRegressor (name unchanged, values updated)
$ws.Range("B3").Value = 0.3949679807973629
$ws.Range("C3").Value = 0.2797384705511541
$ws.Range("D3").Value = -7.343314286554488

# Row 4 - GradientBoostingRegressor -> DecisionTreeRegressor
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.844714538564235
$ws.Range("C4").Value = 0.8463794093002422
$ws.Range("D4").Value = -4.179772337997509

# Row 5 - AdaBoostRegressor -> MLPRegressor
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.2612139991218055
$ws.Range("C5").Value = -0.3146934252327823
$ws.Range("D5").Value = -26.05227733002145
